# CTECH-2652 fix dates as strings (for cut labels)
# Rework the holdings extract columns: drop sub_holding_keys, SourcePortfolioId,
# SourcePortfolioScope and cost_portfolio_ccy.amount columns, and rename the
# remaining headers to the new LUSID extract field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unneeded columns, right-to-left so earlier column letters
# stay valid for subsequent deletes.
$ws.Range("L:L").Delete()   # cost_portfolio_ccy.amount
$ws.Range("F:F").Delete()   # SourcePortfolioScope(default-Properties)
$ws.Range("E:E").Delete()   # SourcePortfolioId(default-Properties)
$ws.Range("C:C").Delete()   # sub_holding_keys

# Rename the remaining headers to match the new extract field names.
$ws.Range("B1").Value = "luid"
$ws.Range("C1").Value = "instrumentName"
$ws.Range("D1").Value = "holdingType"
$ws.Range("E1").Value = "units"
$ws.Range("F1").Value = "settledUnits"
$ws.Range("G1").Value = "costAmount"
$ws.Range("H1").Value = "costCurrency"
$ws.Range("I1").Value = "portfolioCurrency"

# The GBP currency holding row's instrument name is now simply "GBP"
# (previously this column held the "Name" property for CCY_GBP).
$ws.Range("C7").Value = "GBP"
